$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update C52: "Rédaction de rapport" -> split into two distinct strings.
#    We first rename the existing shared string (sole referrer = C52) to the
#    "(design)" variant, then create C53 referencing that same text so the
#    string is shared, and finally change C52 to the final, different text
#    (this appends a brand-new shared string since the old one now has 2
#    referrers).
# ---------------------------------------------------------------------------
$ws.Range("C52").Value = "Rédaction de rapport (design)"

# ---------------------------------------------------------------------------
# 2. Add the new row 53 (continuation of the 2023-09-07 entry is done, this
#    is the new day 2023-09-08 / serial 45177) and copy formatting from the
#    analogous block above (rows 50-52) so number formats / alignment match.
# ---------------------------------------------------------------------------
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A53:A55").PasteSpecial(-4122) | Out-Null

$ws.Range("B50").Copy() | Out-Null
$ws.Range("B53:B55").PasteSpecial(-4122) | Out-Null

$ws.Range("C52").Copy() | Out-Null
$ws.Range("C53:C55").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the values for the new rows.
# ---------------------------------------------------------------------------
$ws.Range("A53").Value = 45177
$ws.Range("B53").Value = "08:30 - 11:30"
$ws.Range("C53").Value = "Rédaction de rapport (design)"

$ws.Range("C52").Value = "Rédaction de rapport (page de titre + design)"

$ws.Range("B54").Value = "11:30 - 12:30"
$ws.Range("C54").Value = "Séance avec maitre de diplôme"

$ws.Range("B55").Value = "13:30 - 16:00"
$ws.Range("C55").Value = "Programmation module RFID"

# ---------------------------------------------------------------------------
# 4. Merge the date column for the new day, like the other day blocks.
# ---------------------------------------------------------------------------
$ws.Range("A53:A55").Merge()

# ---------------------------------------------------------------------------
# 5. Update the view: scroll down and select C60 (as in the authored file).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("C60").Select()

Write-Host "edit complete"
